# The underlying edit (per the XML diff) swaps the contents of
# ppt/theme/theme1.xml ("Office Theme" / default Office colors, used by
# the Notes Master) and ppt/theme/theme2.xml ("Integral" / "Red Violet"
# colors, used by the Slide Master and therefore by every slide).
#
# The PowerPoint object model doesn't give us raw part-swapping, but the
# actually-visible effect of that swap is that the deck's live color
# scheme (driven by theme2.xml, since that's what the Slide Master and
# all slides use) changes from the pink/purple "Red Violet" palette to
# the plain default "Office" palette. We reproduce that with the Theme
# Color Scheme API, writing each of the 12 theme colors (dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink) to the standard Office values that
# theme1.xml already uses.

$p = $ppt.ActivePresentation

# Target colors, in MsoThemeColorSchemeIndex order (1-12), as plain
# RRGGBB hex strings (the values theme1.xml/"Office Theme" already has).
$officeColorsRgbHex = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $hex = $officeColorsRgbHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)

    # Ole/VBA RGB() long packs as 0x00BBGGRR.
    $oleRgb = ($b * 65536) + ($g * 256) + $r

    $colorScheme.Colors($i).RGB = $oleRgb
}

Write-Output "Slide Master theme colors updated to Office defaults."
